# Rever_DailyTrack_BALRAJ_2022.xlsx - "Add files via upload"
#
# Adds a new daily-tracking entry (row 33) to the FEB-22 sheet with a
# second comment line (row 34), matching the formatting already used by
# the existing entries (copied from row 31, which has the same "day
# entry + continuation" pattern as rows 33/34).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEB-22")

# Row 33 is a brand new dated entry: copy the cell formatting (borders,
# shading, number formats, wrap text, ...) from row 31 - an existing
# entry with the exact same column layout - onto rows 33:34 first, then
# fill in the values.
$ws.Range("A31:G31").Copy()
$ws.Range("A33:G33").PasteSpecial(-4122)

$ws.Range("A33").Value = 22
$ws.Range("B33").Value = 44620
$ws.Range("C33").Value = "RPA GSS"
$ws.Range("D33").Value = "1. The task of invoice generation has been completed, tested and it is running smoothly, whereas during testing, mohan san suggesting that after creating csv file, we need to update the master file to avoid memory craking  (now it is done dynamically getting data and updating master file dynamically) and it is`nwork in progress"
$ws.Range("E33").Value = 0.8
$ws.Range("F33").Value = "WIP"
$ws.Rows.Item(33).RowHeight = 57.6

# Row 34 holds the second comment line for the same entry.
$ws.Range("D34").Value = "2. Athough the captching is working fine, still Mohan has suggested to check 30 or 50 times of the captcha intead of continues checking  and we need to enhance the captcha "
$ws.Rows.Item(34).RowHeight = 28.8

$ws.Range("D34").Select()
